# Generate Report for Handoff
# Updates the localization-status report to reflect that a new handoff
# package was generated for e2e\b.md (the en-US source moved on, so the
# previously handed-back translation is now stale and a fresh handoff
# round has started).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ac44cc1da437f8f0c1b6bbb6a98543a11b3fe654/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2ffa087926ce3e3f82cb7592c313ccf375b22d05/e2e/b.md."

# ---- Overview sheet: b.md row (row 3) now shows "Ready for handoff" ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-10-14 07:47:40"

# ---- zh-cn sheet: b.md row (row 3) ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-10-14 07:47:30"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---- de-de sheet: b.md row (row 3) ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-10-14 07:47:40"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
